# Auto-generated Excel COM-interop script
# Applies the cell-value updates described by the commit diff
# across sheets ALC, ARM, CRP, CUL, GSM, LTW, WVR (no changes to BSM).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 40003136
$ws.Range("I62").Value = 66668520
$ws.Range("J62").Value = 5061
$ws.Range("K62").Value = 66668520
$ws.Range("L62").Value = 5061
$ws.Range("M62").Value = -66667896
$ws.Range("N62").Value = -6309
$ws.Range("H64").Value = 3060.3225
$ws.Range("I64").Value = 2698.0952
$ws.Range("K64").Value = 2698.0952
$ws.Range("M64").Value = -2450.0952
$ws.Range("H65").Value = 40003136
$ws.Range("I65").Value = 66668520
$ws.Range("J65").Value = 5061
$ws.Range("K65").Value = 333342600
$ws.Range("L65").Value = 25305
$ws.Range("M65").Value = -333339480
$ws.Range("N65").Value = -31545
$ws.Range("H67").Value = 3060.3225
$ws.Range("I67").Value = 2698.0952
$ws.Range("K67").Value = 2698.0952
$ws.Range("M67").Value = -1840.0952
$ws.Range("H98").Value = 2649.36
$ws.Range("I98").Value = 1300.6
$ws.Range("J98").Value = 8044.4
$ws.Range("K98").Value = 1300.6
$ws.Range("L98").Value = 8044.4
$ws.Range("M98").Value = 197.4000000000001
$ws.Range("N98").Value = -11040.4
$ws.Range("H122").Value = 2649.36
$ws.Range("I122").Value = 1300.6
$ws.Range("J122").Value = 8044.4
$ws.Range("K122").Value = 3901.8
$ws.Range("L122").Value = 24133.2
$ws.Range("M122").Value = -1451.8
$ws.Range("N122").Value = -29033.2
$ws.Range("H125").Value = 1477
$ws.Range("I125").Value = 1310
$ws.Range("J125").Value = 1577.2
$ws.Range("K125").Value = 11790
$ws.Range("L125").Value = 14194.8
$ws.Range("M125").Value = -9330
$ws.Range("N125").Value = -19114.8
$ws.Range("H137").Value = 1070177.9
$ws.Range("I137").Value = 1375.3334
$ws.Range("J137").Value = 3474983.8
$ws.Range("K137").Value = 4126.0002
$ws.Range("L137").Value = 10424951.4
$ws.Range("M137").Value = -1576.0002
$ws.Range("N137").Value = -10430051.4
$ws.Range("H138").Value = 3981.6978
$ws.Range("I138").Value = 3216.5
$ws.Range("J138").Value = 4948.263
$ws.Range("K138").Value = 9649.5
$ws.Range("L138").Value = 14844.789
$ws.Range("M138").Value = -4509.5
$ws.Range("N138").Value = -25124.789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21314.758
$ws.Range("I32").Value = 9051.514999999999
$ws.Range("J32").Value = 33578
$ws.Range("K32").Value = 9051.514999999999
$ws.Range("L32").Value = 33578
$ws.Range("M32").Value = -8764.514999999999
$ws.Range("N32").Value = -34152
$ws.Range("H74").Value = 80199.92999999999
$ws.Range("I74").Value = 100787.4
$ws.Range("K74").Value = 100787.4
$ws.Range("M74").Value = -99913.39999999999
$ws.Range("H77").Value = 80199.92999999999
$ws.Range("I77").Value = 100787.4
$ws.Range("K77").Value = 503937
$ws.Range("M77").Value = -499569
$ws.Range("H102").Value = 1067.8572
$ws.Range("I102").Value = 1067.8572
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1067.8572
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 554.1428000000001
$ws.Range("N102").ClearContents()
$ws.Range("H108").Value = 26000
$ws.Range("J108").Value = 26000
$ws.Range("L108").Value = 26000
$ws.Range("N108").Value = -33680
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 21000
$ws.Range("J112").Value = 21000
$ws.Range("L112").Value = 21000
$ws.Range("N112").Value = -23954
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H117").Value = 29430
$ws.Range("J117").Value = 29430
$ws.Range("L117").Value = 29430
$ws.Range("M117").Value = -38608
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H119").Value = 31200
$ws.Range("J119").Value = 31200
$ws.Range("L119").Value = 31200
$ws.Range("N119").Value = -40876

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 77309
$ws.Range("J141").Value = 86919.164
$ws.Range("L141").Value = 86919.164
$ws.Range("N141").Value = -97279.164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 25140
$ws.Range("J32").Value = 25140
$ws.Range("L32").Value = 75420
$ws.Range("N32").Value = -75986
$ws.Range("H107").Value = 877657.3
$ws.Range("J107").Value = 1010560.3
$ws.Range("L107").Value = 3031680.9
$ws.Range("N107").Value = -3035520.9
$ws.Range("H129").Value = 2244.2778
$ws.Range("I129").Value = 1662.75
$ws.Range("J129").Value = 2709.5
$ws.Range("K129").Value = 4988.25
$ws.Range("L129").Value = 8128.5
$ws.Range("M129").Value = 11.75
$ws.Range("N129").Value = -18128.5
$ws.Range("H131").Value = 2034.4117
$ws.Range("I131").Value = 10000
$ws.Range("J131").Value = 1536.5625
$ws.Range("K131").Value = 30000
$ws.Range("L131").Value = 4609.6875
$ws.Range("M131").Value = -24960
$ws.Range("N131").Value = -14689.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = 29
$ws.Range("H27").Value = 39000
$ws.Range("I27").Value = 39000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 39000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -38834
$ws.Range("N27").ClearContents()
$ws.Range("H102").Value = 4274627.5
$ws.Range("I102").Value = 5556425
$ws.Range("J102").Value = 1969
$ws.Range("K102").Value = 5556425
$ws.Range("L102").Value = 1969
$ws.Range("M102").Value = -5554803
$ws.Range("N102").Value = -5213
$ws.Range("H126").Value = 47620536
$ws.Range("I126").Value = 55556956
$ws.Range("K126").Value = 166670868
$ws.Range("M126").Value = -166668398

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1873.1111
$ws.Range("I7").Value = 1311.6
$ws.Range("J7").Value = 2575
$ws.Range("K7").Value = 1311.6
$ws.Range("L7").Value = 2575
$ws.Range("M7").Value = -1199.6
$ws.Range("N7").Value = -2799
$ws.Range("H40").Value = 4211.4287
$ws.Range("I40").Value = 4142
$ws.Range("J40").Value = 4280.857
$ws.Range("K40").Value = 4142
$ws.Range("L40").Value = 4280.857
$ws.Range("M40").Value = -4006
$ws.Range("N40").Value = -4552.857
$ws.Range("H126").Value = 1873.1111
$ws.Range("I126").Value = 1311.6
$ws.Range("J126").Value = 2575
$ws.Range("K126").Value = 3934.8
$ws.Range("L126").Value = 7725
$ws.Range("M126").Value = -1464.8
$ws.Range("N126").Value = -12665
$ws.Range("H136").Value = 6735.45
$ws.Range("I136").Value = 1333.9333
$ws.Range("J136").Value = 22940
$ws.Range("K136").Value = 4001.7999
$ws.Range("L136").Value = 68820
$ws.Range("M136").Value = -1451.7999
$ws.Range("N136").Value = -73920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 6766.6665
$ws.Range("J30").Value = 6766.6665
$ws.Range("L30").Value = 6766.6665
$ws.Range("N30").Value = -6980.6665
$ws.Range("H133").Value = 47115
$ws.Range("J133").Value = 47115
$ws.Range("L133").Value = 47115
$ws.Range("N133").Value = -57235
$ws.Range("H136").Value = 4373.0967
$ws.Range("I136").Value = 790.8077
$ws.Range("J136").Value = 23001
$ws.Range("K136").Value = 2372.4231
$ws.Range("L136").Value = 69003
$ws.Range("M136").Value = 177.5769
$ws.Range("N136").Value = -74103

Write-Host "Applied all cell updates."
